$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New predicted values for columns B and C, rows 2-67,
# produced by re-running the prediction loop (see commit message).
$data = @(
    @(-0.2909141766055221, -0.2757346137170916),
    @(0.250429317339002, 0.2503386379796715),
    @(0.1809838537612315, 0.1898448635988708),
    @(-0.0356568537171813, -0.03069743402202884),
    @(0.1826398358552733, 0.189986233938121),
    @(-0.4282759348613254, -0.4188081259184758),
    @(-0.2750526674036062, -0.2510779268284625),
    @(-0.3717861371237773, -0.3701072376973022),
    @(0.3472983702021971, 0.3587265422777822),
    @(-0.1952196377481137, -0.1880564208480707),
    @(-0.04841202328792223, -0.0369409108347238),
    @(0.006802512295449493, 0.005706268408748546),
    @(0.1165728043122438, 0.1249572479757306),
    @(0.003090940948384968, 0.01584046856950461),
    @(0.379763138506495, 0.3979898313796012),
    @(0.660010392767101, 0.654067064061202),
    @(0.0831660780375153, 0.05046369703599946),
    @(0.43642760268563, 0.4290197685502667),
    @(0.3415487805647693, 0.3547297180093819),
    @(0.5126516164589306, 0.5479638996118918),
    @(0.424547355049718, 0.426269999648244),
    @(-0.01126124077300952, -0.005631517867795921),
    @(4.334097336613688, 4.317720391123286),
    @(0.4934852282346514, 0.4668881953821007),
    @(0.4275772796777446, 0.4046308426275156),
    @(0.2847042344339095, 0.2586029384636029),
    @(1.04229764095734, 1.002813007074137),
    @(5.662181268714689, 5.170629289836413),
    @(0.9653505985672939, 0.924117484475707),
    @(-0.2184845195128625, -0.2632749190854577),
    @(0.7846543751460326, 0.7495677028461432),
    @(0.8737302907557974, 0.855501533467473),
    @(-0.6682857482188554, -0.6844880111018772),
    @(0.8113291412439705, 0.8099285809726998),
    @(0.7530154080943274, 0.7509634599418302),
    @(0.7334962419193686, 0.7275065559500377),
    @(0.7422586571568021, 0.7251387680663887),
    @(0.564315474342054, 0.5702895628981168),
    @(0.7337100851847743, 0.7413699042840771),
    @(0.5661099619036586, 0.5624140041041261),
    @(0.6848554698128246, 0.6703943901995139),
    @(0.7175126490592517, 0.7083885862451348),
    @(0.6565737711511859, 0.6579767609369422),
    @(0.6153639377726098, 0.6191152717822307),
    @(-1.277621360328754, -1.270880536667543),
    @(-0.99184933780554, -0.9871065606551556),
    @(-0.887860921688258, -0.8816908956008629),
    @(-0.6502216361091679, -0.6441782177273736),
    @(-0.05898409078580606, -0.05502435637365417),
    @(-0.8764437831433545, -0.8678832629530725),
    @(-0.8764437831433545, -0.8678832629530725),
    @(-1.109071977712408, -1.110153319826916),
    @(-0.1989048519688044, -0.1908628337365867),
    @(-1.01369646402942, -1.006858221511667),
    @(-0.9102532362208542, -0.8966458885303238),
    @(-0.9569878829587767, -0.9436111128636898),
    @(-1.155999245559128, -1.129488322093305),
    @(-0.8755454492962101, -0.8533938298445864),
    @(-0.5123768353802537, -0.4888985070244863),
    @(0.366436512377399, 0.3702423988900625),
    @(-1.247693151957668, -1.227488952603415),
    @(-0.758098723716701, -0.7198624288737678),
    @(-0.8966681281711343, -0.8878765045190393),
    @(-0.1328292803076666, -0.1086031543723337),
    @(-0.8236383263975512, -0.7976715977463904),
    @(-0.8255639498988997, -0.7891681206752141)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $data[$i]
    $ws.Cells.Item($i + 2, 2).Value = $row[0]
    $ws.Cells.Item($i + 2, 3).Value = $row[1]
}
